# === Parts list update: insert column, add rows 13-22, restructure hyperlinks ===
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks up-front; they will be re-added at their final locations
# once the new column has been inserted and all data is in place.
$ws.Hyperlinks.Delete()

# Insert a new blank column at F ('No. required'), shifting the old F-I columns to G-J
$ws.Columns("F").Insert()

# ---- Row 4 header: add new header cells for columns E (size), J (Mnfctr), K (Mnfctr part no.), M (Note) ----
$ws.Cells.Item(4, 5).Value = 'purpose'
$ws.Cells.Item(4, 10).Value = 'Mnfctr'
$ws.Cells.Item(4, 11).Value = 'Mnfctr. part no.'
$ws.Cells.Item(4, 13).Value = 'Note'
$ws.Cells.Item(4, 5).Font.Bold = $true
$ws.Cells.Item(4, 10).Font.Bold = $true
$ws.Cells.Item(4, 11).Font.Bold = $true
$ws.Cells.Item(4, 13).Font.Bold = $true

# ---- Row 3 filler cells for newly added columns (J, K, L are blank/styled) ----
$ws.Range("J3:L3").Value = ""

# ---- Populate data rows 5-22 ----
# -- Row 5 --
$ws.Cells.Item(5, 1).Value = 'banana connector'
$ws.Cells.Item(5, 2).Value = 'plug '
$ws.Cells.Item(5, 3).Value = 'black'
$ws.Cells.Item(5, 4).Value = '4mm'
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.5
$ws.Cells.Item(5, 8).Value = 'farnel'
$ws.Cells.Item(5, 9).Value = 1698951
$ws.Cells.Item(5, 10).Value = 'Multicomp'
$ws.Cells.Item(5, 11).Value = '24.247.2'

# -- Row 6 --
$ws.Cells.Item(6, 1).Value = 'banana connector'
$ws.Cells.Item(6, 2).Value = 'plug'
$ws.Cells.Item(6, 3).Value = 'red'
$ws.Cells.Item(6, 4).Value = '4mm'
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.5
$ws.Cells.Item(6, 8).Value = 'farnel'
$ws.Cells.Item(6, 9).Value = 1698950
$ws.Cells.Item(6, 10).Value = 'Multicomp'
$ws.Cells.Item(6, 11).Value = '24.247.1'

# -- Row 7 --
$ws.Cells.Item(7, 1).Value = 'banana connector'
$ws.Cells.Item(7, 2).Value = 'jack'
$ws.Cells.Item(7, 3).Value = 'black'
$ws.Cells.Item(7, 4).Value = '4mm'
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.35
$ws.Cells.Item(7, 8).Value = 'farnel'
$ws.Cells.Item(7, 9).Value = 1698964
$ws.Cells.Item(7, 10).Value = 'Multicomp'
$ws.Cells.Item(7, 11).Value = '25.413.2'

# -- Row 8 --
$ws.Cells.Item(8, 1).Value = 'banana connector'
$ws.Cells.Item(8, 2).Value = 'jack'
$ws.Cells.Item(8, 3).Value = 'red'
$ws.Cells.Item(8, 4).Value = '4mm'
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.35
$ws.Cells.Item(8, 8).Value = 'farnel'
$ws.Cells.Item(8, 9).Value = 1698963
$ws.Cells.Item(8, 10).Value = 'Multicomp'
$ws.Cells.Item(8, 11).Value = '25.413.1'

# -- Row 9 --
$ws.Cells.Item(9, 1).Value = 'DC power connector'
$ws.Cells.Item(9, 2).Value = 'plug'
$ws.Cells.Item(9, 3).Value = 'cable mount'
$ws.Cells.Item(9, 4).Value = '2.5mm'
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.5
$ws.Cells.Item(9, 8).Value = 'farnel'
$ws.Cells.Item(9, 9).Value = 1737249
$ws.Cells.Item(9, 10).Value = 'Multicomp'
$ws.Cells.Item(9, 11).Value = 'MP-122M'
$ws.Cells.Item(9, 13).Value = 'not needed if power supply used'

# -- Row 10 --
$ws.Cells.Item(10, 1).Value = 'DC power connector'
$ws.Cells.Item(10, 2).Value = 'jack'
$ws.Cells.Item(10, 3).Value = 'panel mount'
$ws.Cells.Item(10, 4).Value = '2.5mm'
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9
$ws.Cells.Item(10, 8).Value = 'farnel'
$ws.Cells.Item(10, 9).Value = 1737252
$ws.Cells.Item(10, 10).Value = 'Multicomp'
$ws.Cells.Item(10, 11).Value = 'MJ-15SR'

# -- Row 11 --
$ws.Cells.Item(11, 1).Value = '2-pole audio connector'
$ws.Cells.Item(11, 2).Value = 'plug'
$ws.Cells.Item(11, 3).Value = 'cable mount'
$ws.Cells.Item(11, 4).Value = '2.5mm'
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.73
$ws.Cells.Item(11, 8).Value = 'farnel'
$ws.Cells.Item(11, 9).Value = 1243261
$ws.Cells.Item(11, 10).Value = 'Lumberg'
$ws.Cells.Item(11, 11).Value = 'KLS 2'

# -- Row 12 --
$ws.Cells.Item(12, 1).Value = '2-pole audio connector'
$ws.Cells.Item(12, 2).Value = 'jack'
$ws.Cells.Item(12, 3).Value = 'panel mount'
$ws.Cells.Item(12, 4).Value = '2.5mm'
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.24
$ws.Cells.Item(12, 8).Value = 'farnel'
$ws.Cells.Item(12, 9).Value = 1267394
$ws.Cells.Item(12, 10).Value = 'Pro Signal'
$ws.Cells.Item(12, 11).Value = 'MJ-164H'

# -- Row 13 --
$ws.Cells.Item(13, 1).Value = 'DC power supply'
$ws.Cells.Item(13, 2).Value = 'mains AC adaptor'
$ws.Cells.Item(13, 3).Value = '12V, 10W output'
$ws.Cells.Item(13, 4).Value = '2.5mm jack'
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 7.5
$ws.Cells.Item(13, 8).Value = 'farnel'
$ws.Cells.Item(13, 9).Value = 1827445
$ws.Cells.Item(13, 10).Value = 'Multicomp'
$ws.Cells.Item(13, 11).Value = 'MCPLG12V10WUK'

# -- Row 14 --
$ws.Cells.Item(14, 1).Value = 'On-off (rocker) switch'
$ws.Cells.Item(14, 2).Value = 'SPST'
$ws.Cells.Item(14, 3).Value = '250V, 16A max, illuminated'
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1.16
$ws.Cells.Item(14, 8).Value = 'farnel'
$ws.Cells.Item(14, 9).Value = 1454386
$ws.Cells.Item(14, 10).Value = 'Multicomp'
$ws.Cells.Item(14, 11).Value = 'MC34231-091-72'

# -- Row 15 --
$ws.Cells.Item(15, 1).Value = 'On-off (rocker) switch'
$ws.Cells.Item(15, 2).Value = 'SPST'
$ws.Cells.Item(15, 3).Value = '20A, non-illuminated, black/red'
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.64
$ws.Cells.Item(15, 8).Value = 'farnel'
$ws.Cells.Item(15, 9).Value = 1454382
$ws.Cells.Item(15, 10).Value = 'Multicomp'
$ws.Cells.Item(15, 11).Value = 'MC34224-071-1501'

# -- Row 16 --
$ws.Cells.Item(16, 1).Value = 'Potentiometer'
$ws.Cells.Item(16, 2).Value = 'Rotor'
$ws.Cells.Item(16, 3).Value = '10k, 200mW, +-20%'
$ws.Cells.Item(16, 5).Value = 'RG'
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.61
$ws.Cells.Item(16, 8).Value = 'farnel'
$ws.Cells.Item(16, 9).Value = 1760793
$ws.Cells.Item(16, 10).Value = 'BI Technologies'
$ws.Cells.Item(16, 11).Value = 'P160KNP-0QC20B10K'
$ws.Cells.Item(16, 13).Value = 'this or 249BBHS0XB25103KA'

# -- Row 17 --
$ws.Cells.Item(17, 1).Value = 'Potentiometer'
$ws.Cells.Item(17, 2).Value = 'Rotor'
$ws.Cells.Item(17, 3).Value = '10k, 1W, +-10%'
$ws.Cells.Item(17, 5).Value = 'RG'
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 6.94
$ws.Cells.Item(17, 8).Value = 'farnel'
$ws.Cells.Item(17, 9).Value = 9609245
$ws.Cells.Item(17, 10).Value = 'Vishay'
$ws.Cells.Item(17, 11).Value = '249BBHS0XB25103KA'
$ws.Cells.Item(17, 13).Value = 'this or P160KNP-0QC20B10K'

# -- Row 18 --
$ws.Cells.Item(18, 1).Value = 'BNC panel mount'
$ws.Cells.Item(18, 2).Value = 'jack'
$ws.Cells.Item(18, 3).Value = 'brass, <11GHz'
$ws.Cells.Item(18, 4).Value = '50 Ohm'
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 1.67
$ws.Cells.Item(18, 8).Value = 'farnel'
$ws.Cells.Item(18, 9).Value = 1169699
$ws.Cells.Item(18, 10).Value = 'Multicomp'
$ws.Cells.Item(18, 11).Value = '13-25'

# -- Row 19 --
$ws.Cells.Item(19, 1).Value = 'Test point'
$ws.Cells.Item(19, 2).Value = 'thru-hole'
$ws.Cells.Item(19, 3).Value = 'for panel mounting'
$ws.Cells.Item(19, 4).Value = '1mm diamter'
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 1.41
$ws.Cells.Item(19, 8).Value = 'farnel'
$ws.Cells.Item(19, 9).Value = 1702006
$ws.Cells.Item(19, 10).Value = 'Multicomp'
$ws.Cells.Item(19, 11).Value = 'TEST-3'
$ws.Cells.Item(19, 13).Value = 'pack of 100, alt. to 1035'

# -- Row 20 --
$ws.Cells.Item(20, 1).Value = 'Test point'
$ws.Cells.Item(20, 2).Value = 'thru-hole'
$ws.Cells.Item(20, 3).Value = 'for panel mounting'
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 13.7
$ws.Cells.Item(20, 8).Value = 'farnel'
$ws.Cells.Item(20, 9).Value = 2112492
$ws.Cells.Item(20, 10).Value = 'Keystone'
$ws.Cells.Item(20, 11).Value = 1035
$ws.Cells.Item(20, 11).NumberFormat = "@"
$ws.Cells.Item(20, 13).Value = 'pack of 100, alt. to TEST-3'

# -- Row 21 --
$ws.Cells.Item(21, 1).Value = 'Potentiometer'
$ws.Cells.Item(21, 2).Value = 'thru-hole'
$ws.Cells.Item(21, 3).Value = '5k, 1W, +-10%'
$ws.Cells.Item(21, 5).Value = 'V_offset'
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 6.33
$ws.Cells.Item(21, 8).Value = 'farnel'
$ws.Cells.Item(21, 9).Value = 9609237
$ws.Cells.Item(21, 10).Value = 'Vishay'
$ws.Cells.Item(21, 11).Value = '249BBHS0XB25502KA'

# -- Row 22 --
$ws.Cells.Item(22, 1).Value = 'Switch'
$ws.Cells.Item(22, 2).Value = 'SPDT'
$ws.Cells.Item(22, 3).Value = 'panel mount, On-On'
$ws.Cells.Item(22, 4).Value = 'max 5A'
$ws.Cells.Item(22, 5).Value = 'RG_select'
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 1.4
$ws.Cells.Item(22, 8).Value = 'farnel'
$ws.Cells.Item(22, 9).Value = 9473378
$ws.Cells.Item(22, 10).Value = 'Multicomp'
$ws.Cells.Item(22, 11).Value = '1MS1T1B5M1QE'

# ---- L column: Net cost formula (F*G), shared across L5:L22 ----
$ws.Range("L5").Formula = "=F5*G5"
$ws.Range("L6:L22").Formula = "=F6*G6"

# ---- J5:J11 Mnfctr column uses the existing grey font style (matches Supplier order code column) ----
$ws.Range("J5:J11").Font.Color = 3355443
$ws.Range("J14:J22").Font.Color = 3355443

# ---- Re-create hyperlinks on column H (farnel) at their final positions ----
$ws.Hyperlinks.Add($ws.Cells.Item(5, 8), 'http://uk.farnell.com/multicomp/24-247-2/receptacle-32a-4mm-panel-black/dp/1698951') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 8), 'http://uk.farnell.com/multicomp/24-247-1/receptacle-32a-4mm-panel-red/dp/1698950') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 8), 'http://uk.farnell.com/multicomp/25-413-2/plug-24a-4mm-cable-black/dp/1698964') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 8), 'http://uk.farnell.com/multicomp/25-413-1/banana-plug-24a-4mm-cable-red/dp/1698963') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 8), 'http://uk.farnell.com/multicomp/mp-122m/plug-dc-power-2-5mm/dp/1737249', 'techDocsHook') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 8), 'http://uk.farnell.com/multicomp/mj-15sr/chassis-socket-psu-panel-mount/dp/1737252') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 8), 'http://uk.farnell.com/lumberg/kls-2/connector-plug-3-5mm-mono-cable/dp/1243261') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 8), 'http://uk.farnell.com/lumberg/kls-2/connector-plug-3-5mm-mono-cable/dp/1267394') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 8), 'http://uk.farnell.com/multicomp/mcplg12v10wuk/mains-ac-adaptor-12v-10w/dp/1827445') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14, 8), 'http://uk.farnell.com/multicomp/mc34231-091-72/rocker-switch-spst-illuminated/dp/1454386') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15, 8), 'http://uk.farnell.com/multicomp/mc34224-071-1501/rocker-switch-spst-black-red/dp/1454382') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16, 8), 'http://uk.farnell.com/bi-technologies/p160knp-0qc20b10k/potentiometer-10k/dp/1760793') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17, 8), 'http://uk.farnell.com/vishay/249bbhs0xb25103ka/potentiometer-10k-1w/dp/9609245') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(18, 8), 'http://uk.farnell.com/pro-signal/mj-164h/bnc-panel-socket/dp/1169699') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(19, 8), 'http://uk.farnell.com/keystone/test-3/test-point-red/dp/1702006') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(20, 8), 'http://uk.farnell.com/keystone/1035/test-point-pack-100/dp/2112492') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(21, 8), 'http://uk.farnell.com/vishay/249bbhs0xb25502ka/potentiometer-5k-1w/dp/9609237') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(22, 8), 'http://uk.farnell.com/multicomp/1ms1t1b5m1qe/toggle-switch-spdt/dp/9473378') | Out-Null

# ---- Column widths for new columns E/F/J/K (approximate the target layout) ----
$ws.Columns("F").ColumnWidth = 6.6640625
$ws.Columns("J").ColumnWidth = 17.6640625
$ws.Columns("K").ColumnWidth = 17.6640625

# ---- View settings: match target top-left cell & active selection ----
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("M10").Select()
